$d = $word.ActiveDocument

# 1) Title text change: "My Personal Website" -> "My way of expressing" inside <title>...</title>
$d.Content.Find.Execute("<title>My Personal Website</title>", $true, $true, $false, $false, $false, $true, 1, $false, "<title>My way of expressing</title>", 2) | Out-Null

# 2) Replace the logo block (<div class="logo">, <img .../>, </div>) and the old <h1> text
#    with a single new <h1> paragraph: "Welcome lets make it fun".
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($startPara -eq $null -and $txt -match '<div class="logo">') {
        $startPara = $i
    }
    if ($txt -match '<h1>Welcome to My Personal Website</h1>') {
        $endPara = $i
        break
    }
}

# Delete the 3 paragraphs preceding the <h1> paragraph (div, img, /div),
# leaving the <h1> paragraph itself in place so we can rewrite its text.
$delStart = $d.Paragraphs.Item($startPara).Range.Start
$delEnd = $d.Paragraphs.Item($endPara - 1).Range.End
$d.Range($delStart, $delEnd).Delete() | Out-Null

$headingPara = $d.Paragraphs.Item($startPara)
$headingPara.Range.Text = "        <h1>Welcome lets make it fun</h1>"

# 3) Footer text change
$d.Content.Find.Execute("&copy; 2025 My Personal Website", $true, $true, $false, $false, $false, $true, 1, $false, "&copy; Let me express", 2) | Out-Null
